$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7:29 (E_HAPD ... E_TEUR) down to rows 8:30 to make room for a
# new "E_HA" row. Use the single-call Range.Copy(Destination) form (rather
# than separate Copy() + PasteSpecial()) so cell formatting -- including the
# shared style index "4" used by every data row, and row 30 which lies
# beyond the original A1:C29 dimension -- is carried over faithfully instead
# of a new style being fabricated or dropped.
$ws.Range("A7:C29").Copy($ws.Range("A8:C30"))

# Row 20 (E_NN) had blank Einheit De / Einheit En cells; after the shift it
# is now row 21. Copying a blank source cell does not clear the destination
# cell's previous contents, so clear these explicitly.
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = ""

# Fill in the newly inserted row 7 with the new unit "E_HA" / "Hektar" / "Hectre".
$ws.Range("A7").Value = "E_HA"
$ws.Range("B7").Value = "Hektar"
$ws.Range("C7").Value = "Hectre"
